$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 9 (shifts subsequent rows, e.g. former row 10, up by one)
$ws.Rows.Item(9).Delete()

# Update the active sheet view: scroll so row 4 is the top-left visible row,
# and select the entire (new) row 9
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A9:XFD9").Select()
